# The workbook author cleared out the roster's First/Last Name columns
# (the "i dont even remember the changes" commit) for all four
# division tables on the Roster sheet, leaving the Weight column intact.
# Everything else in the workbook (Leaderboards, Points calculator,
# Juniors, Division 1/2/3) just recalculates off of that change.

$wb = $excel.ActiveWorkbook

$roster = $wb.Worksheets.Item("Roster")

# roster_division1 (Pro 1): First Name / Last Name columns A:B, rows 4-15
$roster.Range("A4:B15").ClearContents()
# roster_division2 (Pro 2): First Name / Last Name columns E:F, rows 4-15
$roster.Range("E4:F15").ClearContents()
# roster_division3 (Pro 3): First Name / Last Name columns I:J, rows 4-15
$roster.Range("I4:J15").ClearContents()
# roster_division0 (Junior): First Name / Last Name columns M:N, rows 4-15
$roster.Range("M4:N15").ClearContents()

# Selection left on the Roster sheet after the edit
[void]$roster.Range("B18").Select()

# Recalculate everything that depends on the roster names
$excel.CalculateFull()

# Move focus to the "Points calculator" sheet, which ends up the
# active tab / selected sheet after the edit
$pc = $wb.Worksheets.Item("Points calculator")
[void]$pc.Activate()
[void]$pc.Range("E8").Select()
